$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1.596690938351973
$ws.Range("F2").Value = 3.124404084168053
$ws.Range("E3").Value = 1.596690938351973
$ws.Range("F3").Value = 3.124404084168053
$ws.Range("E4").Value = 1.596690938351973
$ws.Range("F4").Value = 3.124404084168053
$ws.Range("E5").Value = 1.596690938351973
$ws.Range("F5").Value = 3.124404084168053
$ws.Range("E6").Value = 3.621782456732601
$ws.Range("F6").Value = 1.058594744827242
$ws.Range("E7").Value = 3.621782456732601
$ws.Range("F7").Value = 1.058594744827242
$ws.Range("E8").Value = 3.621782456732601
$ws.Range("F8").Value = 1.058594744827242
$ws.Range("E9").Value = 3.621782456732601
$ws.Range("F9").Value = 1.058594744827242
$ws.Range("E10").Value = 2.609236697542287
$ws.Range("F10").Value = 2.091499414497647
$ws.Range("E11").Value = 2.609236697542287
$ws.Range("F11").Value = 2.091499414497647
$ws.Range("E12").Value = 2.609236697542287
$ws.Range("F12").Value = 2.091499414497647
$ws.Range("E13").Value = 2.609236697542287
$ws.Range("F13").Value = 2.091499414497647
$ws.Range("E14").Value = 2.609236697542287
$ws.Range("F14").Value = 2.091499414497647
$ws.Range("E15").Value = 2.609236697542287
$ws.Range("F15").Value = 2.091499414497647
$ws.Range("E16").Value = 5.64687397511323
$ws.Range("F16").Value = 3.124404084168053
$ws.Range("E17").Value = 5.64687397511323
$ws.Range("F17").Value = 3.124404084168053
$ws.Range("E18").Value = 5.64687397511323
$ws.Range("F18").Value = 3.124404084168053
$ws.Range("E19").Value = 5.64687397511323
$ws.Range("F19").Value = 3.124404084168053
$ws.Range("E20").Value = 7.671965493493858
$ws.Range("F20").Value = 1.058594744827242
$ws.Range("E21").Value = 7.671965493493858
$ws.Range("F21").Value = 1.058594744827242
$ws.Range("E22").Value = 7.671965493493858
$ws.Range("F22").Value = 1.058594744827242
$ws.Range("E23").Value = 7.671965493493858
$ws.Range("F23").Value = 1.058594744827242
$ws.Range("E24").Value = 6.659419734303544
$ws.Range("F24").Value = 2.091499414497647
$ws.Range("E25").Value = 6.659419734303544
$ws.Range("F25").Value = 2.091499414497647
$ws.Range("E26").Value = 6.659419734303544
$ws.Range("F26").Value = 2.091499414497647
$ws.Range("E27").Value = 6.659419734303544
$ws.Range("F27").Value = 2.091499414497647
$ws.Range("E28").Value = 6.659419734303544
$ws.Range("F28").Value = 2.091499414497647
$ws.Range("E29").Value = 6.659419734303544
$ws.Range("F29").Value = 2.091499414497647
$ws.Range("E30").Value = 22.90279553039436
$ws.Range("F30").Value = 2.154094431234792
$ws.Range("E31").Value = 9.29000350637239
$ws.Range("F31").Value = 3.106467756947701
$ws.Range("E32").Value = 9.29000350637239
$ws.Range("F32").Value = 3.106467756947701
$ws.Range("E33").Value = 9.29000350637239
$ws.Range("F33").Value = 3.106467756947701
$ws.Range("E34").Value = 9.29000350637239
$ws.Range("F34").Value = 3.106467756947701
$ws.Range("E35").Value = 11.06326127902205
$ws.Range("F35").Value = 1.078171549326483
$ws.Range("E36").Value = 11.06326127902205
$ws.Range("F36").Value = 1.078171549326483
$ws.Range("E37").Value = 11.06326127902205
$ws.Range("F37").Value = 1.078171549326483
$ws.Range("E38").Value = 11.06326127902205
$ws.Range("F38").Value = 1.078171549326483
$ws.Range("E39").Value = 10.17663239269722
$ws.Range("F39").Value = 2.092319653137092
$ws.Range("E40").Value = 10.17663239269722
$ws.Range("F40").Value = 2.092319653137092
$ws.Range("E41").Value = 10.17663239269722
$ws.Range("F41").Value = 2.092319653137092
$ws.Range("E42").Value = 10.17663239269722
$ws.Range("F42").Value = 2.092319653137092
$ws.Range("E43").Value = 10.17663239269722
$ws.Range("F43").Value = 2.092319653137092
$ws.Range("E44").Value = 10.17663239269722
$ws.Range("F44").Value = 2.092319653137092
$ws.Range("E45").Value = 12.8365190516717
$ws.Range("F45").Value = 3.106467756947701
$ws.Range("E46").Value = 12.8365190516717
$ws.Range("F46").Value = 3.106467756947701
$ws.Range("E47").Value = 12.8365190516717
$ws.Range("F47").Value = 3.106467756947701
$ws.Range("E48").Value = 12.8365190516717
$ws.Range("F48").Value = 3.106467756947701
$ws.Range("E49").Value = 14.60977682432135
$ws.Range("F49").Value = 1.078171549326483
$ws.Range("E50").Value = 14.60977682432135
$ws.Range("F50").Value = 1.078171549326483
$ws.Range("E51").Value = 14.60977682432135
$ws.Range("F51").Value = 1.078171549326483
$ws.Range("E52").Value = 14.60977682432135
$ws.Range("F52").Value = 1.078171549326483
$ws.Range("E53").Value = 13.72314793799652
$ws.Range("F53").Value = 2.092319653137092
$ws.Range("E54").Value = 13.72314793799652
$ws.Range("F54").Value = 2.092319653137092
$ws.Range("E55").Value = 13.72314793799652
$ws.Range("F55").Value = 2.092319653137092
$ws.Range("E56").Value = 13.72314793799652
$ws.Range("F56").Value = 2.092319653137092
$ws.Range("E57").Value = 13.72314793799652
$ws.Range("F57").Value = 2.092319653137092
$ws.Range("E58").Value = 13.72314793799652
$ws.Range("F58").Value = 2.092319653137092
$ws.Range("E59").Value = 16.383034596971
$ws.Range("F59").Value = 3.106467756947701
$ws.Range("E60").Value = 16.383034596971
$ws.Range("F60").Value = 3.106467756947701
$ws.Range("E61").Value = 16.383034596971
$ws.Range("F61").Value = 3.106467756947701
$ws.Range("E62").Value = 16.383034596971
$ws.Range("F62").Value = 3.106467756947701
$ws.Range("E63").Value = 18.15629236962066
$ws.Range("F63").Value = 1.078171549326483
$ws.Range("E64").Value = 18.15629236962066
$ws.Range("F64").Value = 1.078171549326483
$ws.Range("E65").Value = 18.15629236962066
$ws.Range("F65").Value = 1.078171549326483
$ws.Range("E66").Value = 18.15629236962066
$ws.Range("F66").Value = 1.078171549326483
$ws.Range("E67").Value = 17.26966348329583
$ws.Range("F67").Value = 2.092319653137092
$ws.Range("E68").Value = 17.26966348329583
$ws.Range("F68").Value = 2.092319653137092
$ws.Range("E69").Value = 17.26966348329583
$ws.Range("F69").Value = 2.092319653137092
$ws.Range("E70").Value = 17.26966348329583
$ws.Range("F70").Value = 2.092319653137092
$ws.Range("E71").Value = 17.26966348329583
$ws.Range("F71").Value = 2.092319653137092
$ws.Range("E72").Value = 17.26966348329583
$ws.Range("F72").Value = 2.092319653137092
$ws.Range("E73").Value = 0.306165312289297
$ws.Range("F73").Value = 2.260389217068957
$ws.Range("E74").Value = 0.306165312289297
$ws.Range("F74").Value = 2.260389217068957
$ws.Range("E75").Value = 0.306165312289297
$ws.Range("F75").Value = 2.260389217068957
$ws.Range("E76").Value = 0.306165312289297
$ws.Range("F76").Value = 2.260389217068957
